$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on Tue Oct  1 19:11:04 UTC 2024.
# Each cell is forced to Text format before the write and reset to the Normal style
# afterwards, so price strings such as "1.00" / "0.999" are not silently coerced into
# numbers by Excel's COM type inference (matches the workbook's original inlineStr cells).
function Set-TextCell($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '61.919.92'
Set-TextCell 'E2' '  -2.44%  '

# Row 3
Set-TextCell 'D3' '2.501.43'
Set-TextCell 'E3' '  -3.42%  '

# Row 4
Set-TextCell 'E4' '  +0.06%  '

# Row 5
Set-TextCell 'D5' '551.04'
Set-TextCell 'E5' '  -3.71%  '

# Row 6
Set-TextCell 'D6' '147.44'
Set-TextCell 'E6' '  -4.93%  '

# Row 7
Set-TextCell 'E7' '  +0.03%  '

# Row 8
Set-TextCell 'D8' '0.620'
Set-TextCell 'E8' '  -0.48%  '

# Row 9
Set-TextCell 'D9' '2.501.41'
Set-TextCell 'E9' '  -3.35%  '

# Row 10
Set-TextCell 'E10' '  -9.01%  '

# Row 11
Set-TextCell 'E11' '  -1.45%  '

# Row 12
Set-TextCell 'D12' '5.37'
Set-TextCell 'E12' '  -8.21%  '

# Row 13
Set-TextCell 'D13' '0.356'
Set-TextCell 'E13' '  -6.25%  '

# Row 14
Set-TextCell 'D14' '26.16'
Set-TextCell 'E14' '  -7.23%  '

# Row 15
Set-TextCell 'D15' '2.951.31'
Set-TextCell 'E15' '  -3.51%  '

# Row 16
Set-TextCell 'D16' '61.844.10'
Set-TextCell 'E16' '  -2.27%  '

# Row 17
Set-TextCell 'D17' '0.0000164'
Set-TextCell 'E17' '  -8.13%  '

# Row 18
Set-TextCell 'D18' '2.498.13'
Set-TextCell 'E18' '  -2.98%  '

# Row 19
Set-TextCell 'D19' '11.14'
Set-TextCell 'E19' '  -6.94%  '

# Row 20
Set-TextCell 'E20' '  -6.57%  '

# Row 21
Set-TextCell 'D21' '4.19'
Set-TextCell 'E21' '  -7.73%  '

# Row 22
Set-TextCell 'D22' '321.86'
Set-TextCell 'E22' '  -5.96%  '

# Row 23
Set-TextCell 'D23' '1.00'
Set-TextCell 'E23' '  +0.10%  '

# Row 24
Set-TextCell 'D24' '63.87'
Set-TextCell 'E24' '  -5.10%  '

# Row 25
Set-TextCell 'D25' '1.74'
Set-TextCell 'E25' '  -4.46%  '

# Row 26
Set-TextCell 'E26' '  -5.74%  '

# Row 27
Set-TextCell 'D27' '2.625.82'
Set-TextCell 'E27' '  -3.13%  '

# Row 28
Set-TextCell 'B28' 'Fetch.AI'
Set-TextCell 'C28' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D28' '1.49'
Set-TextCell 'E28' '  -4.23%  '

# Row 29
Set-TextCell 'B29' 'Binance-PegBSC-USD'
Set-TextCell 'C29' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell 'D29' '0.999'
Set-TextCell 'E29' '  -0.17%  '

# Row 30
Set-TextCell 'B30' 'Bittensor'
Set-TextCell 'C30' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D30' '538.56'
Set-TextCell 'E30' '  -7.17%  '

# Row 31
Set-TextCell 'D31' '8.37'
Set-TextCell 'E31' '  -8.00%  '

# Row 32
Set-TextCell 'D32' '7.70'
Set-TextCell 'E32' '  -2.26%  '

# Row 33
Set-TextCell 'D33' '0.150'
Set-TextCell 'E33' '  -6.75%  '

# Row 34
Set-TextCell 'D34' '1.90'
Set-TextCell 'E34' '  -7.52%  '

# Row 35
Set-TextCell 'D35' '1.57'
Set-TextCell 'E35' '  -8.81%  '

# Row 36
Set-TextCell 'E36' '  -9.45%  '

# Row 37
Set-TextCell 'D37' '4.88'
Set-TextCell 'E37' '  -10.28%  '

# Row 38
Set-TextCell 'E38' '  +0.01%  '

# Row 39
Set-TextCell 'E39' '  -5.62%  '

# Row 40
Set-TextCell 'D40' '18.56'
Set-TextCell 'E40' '  -5.84%  '

# Row 41
Set-TextCell 'D41' '143.79'
Set-TextCell 'E41' '  -6.86%  '

# Row 42
Set-TextCell 'D42' '0.999'
Set-TextCell 'E42' '  -0.02%  '

# Row 43
Set-TextCell 'E43' '  -8.42%  '

# Row 44
Set-TextCell 'D44' '40.46'
Set-TextCell 'E44' '  -1.93%  '

# Row 45
Set-TextCell 'D45' '2.30'
Set-TextCell 'E45' '  -6.97%  '

# Row 46
Set-TextCell 'D46' '149.24'
Set-TextCell 'E46' '  -4.36%  '

# Row 47
Set-TextCell 'D47' '3.58'
Set-TextCell 'E47' '  -8.46%  '

# Row 48
Set-TextCell 'D48' '20.87'
Set-TextCell 'E48' '  -9.85%  '

# Row 49
Set-TextCell 'D49' '0.0535'
Set-TextCell 'E49' '  -8.87%  '

# Row 50
Set-TextCell 'D50' '0.592'
Set-TextCell 'E50' '  -5.25%  '

# Row 51
Set-TextCell 'D51' '0.0957'
Set-TextCell 'E51' '  -4.70%  '

